# Updated NGIN,VLV & Domain Management code
# Replaces the "16"-suffixed Domain Management sample-data row with a new
# "23"/"24"-suffixed row (DomainMgmt23) on the DomainManagement sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data updates (old "16" sample -> new "23"/"24" sample) ---
$ws.Range("A2").Value  = "DomainMgmt23"
$ws.Range("B2").Value  = "*"
$ws.Range("C2").Value  = "DomainMgmt23"
$ws.Range("D2").Value  = "dmdomain23.com"
$ws.Range("F2").Value  = "domainocn23"
$ws.Range("G2").Value  = "testreference23"
$ws.Range("H2").Value  = "domaincontact23"
$ws.Range("J2").Value  = "domain23@test.com"
$ws.Range("Q2").Value  = "DomainMgmt23"
$ws.Range("R2").Value  = "YES"
$ws.Range("S2").Value  = "NO"
$ws.Range("T2").Value  = "DomainOrder_23"
$ws.Range("U2").Value  = "DomainRFI_23"
$ws.Range("V2").Value  = "DomainOrder_23"
$ws.Range("W2").Value  = "DomainRFI_23"
$ws.Range("Y2").Value  = "DomainService_24"
$ws.Range("AA2").Value = "domainservice_24@gmail.com"
$ws.Range("AD2").Value = "Domainuser23"
$ws.Range("AE2").Value = "domainemail23@gmail.com"
$ws.Range("AJ2").Value = "servicecomp23"
$ws.Range("BH2").Value = "DomainOrderedit_23"
$ws.Range("BI2").Value = "DomainRFIedit_23"
$ws.Range("BJ2").Value = "DomainOrder_23"
$ws.Range("BK2").Value = "DomainRFI_23"

# --- Column width follow-up (new values are longer, columns A & C widen) ---
# Target OOXML width is 15.140625 characters; the COM ColumnWidth setter
# snaps to whole-pixel steps in this runtime, so 14.3 is the closest input
# that lands on the nearest achievable bucket (~15.1667).
$ws.Columns.Item(1).ColumnWidth = 14.3
$ws.Columns.Item(3).ColumnWidth = 14.3

# --- View state: scroll/selection moved from AR1/BF6 to V1/AA6 ---
$ws.Activate()
$ws.Range("AA6").Select()
$excel.ActiveWindow.ScrollColumn = 22
$excel.ActiveWindow.ScrollRow = 1
